# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Masamune_Profits sheets
# as described by the source diff (per-sheet Leve profit tables).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: columns H,I,J,K,L,M,N
$ws.Range("H40").Value = 10464.333
$ws.Range("I40").Value = 21800.4
$ws.Range("J40").Value = 2367.1428
$ws.Range("K40").Value = 21800.4
$ws.Range("L40").Value = 2367.1428
$ws.Range("M40").Value = -21625.4
$ws.Range("N40").Value = -2717.1428

# Row 82: columns H,I,J,K,L,M,N
$ws.Range("H82").Value = 83336970
$ws.Range("I82").Value = 3508.889
$ws.Range("J82").Value = 333337340
$ws.Range("K82").Value = 10526.667
$ws.Range("L82").Value = 1000012020
$ws.Range("M82").Value = -10120.667
$ws.Range("N82").Value = -1000012832

# Row 85: columns H,I,J,K,L,M,N
$ws.Range("H85").Value = 83336970
$ws.Range("I85").Value = 3508.889
$ws.Range("J85").Value = 333337340
$ws.Range("K85").Value = 10526.667
$ws.Range("L85").Value = 1000012020
$ws.Range("M85").Value = -9122.667000000001
$ws.Range("N85").Value = -1000014828

# Row 113: columns H,I,J,K,L,M,N
$ws.Range("H113").Value = 2250
$ws.Range("I113").Value = 2050
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 2050
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 1204
$ws.Range("N113").Value = -8808

# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 15826.627
$ws.Range("I132").Value = 2482.4814
$ws.Range("J132").Value = 71256.16
$ws.Range("K132").Value = 7447.4442
$ws.Range("L132").Value = 213768.48
$ws.Range("M132").Value = -4917.4442
$ws.Range("N132").Value = -218828.48

# Row 138: columns H,I,J,K,L,M,N
$ws.Range("H138").Value = 1566.39
$ws.Range("I138").Value = 814.4651
$ws.Range("J138").Value = 2133.6316
$ws.Range("K138").Value = 2443.3953
$ws.Range("L138").Value = 6400.8948
$ws.Range("M138").Value = 2696.6047
$ws.Range("N138").Value = -16680.8948

$ws = $wb.Worksheets.Item("ARM")
# Row 24: columns H,J,L,N
$ws.Range("H24").Value = 37582.332
$ws.Range("J24").Value = 37582.332
$ws.Range("L24").Value = 37582.332
$ws.Range("N24").Value = -38330.332

# Row 32: columns H,I,K,M
$ws.Range("H32").Value = 11913.164
$ws.Range("I32").Value = 11541.392
$ws.Range("K32").Value = 11541.392
$ws.Range("M32").Value = -11254.392

# Row 61: columns H,I,J,K,L,M,N
$ws.Range("H61").Value = 1415.7567
$ws.Range("I61").Value = 1122.7587
$ws.Range("J61").Value = 2477.875
$ws.Range("K61").Value = 1122.7587
$ws.Range("L61").Value = 2477.875
$ws.Range("M61").Value = -910.7587000000001
$ws.Range("N61").Value = -2901.875

# Row 63: columns H,I,J,K,L,M,N
$ws.Range("H63").Value = 3871.875
$ws.Range("I63").Value = 2829.1667
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 2829.1667
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = -2143.1667
$ws.Range("N63").Value = -8372

# Row 66: columns H,I,J,K,L,M,N
$ws.Range("H66").Value = 3871.875
$ws.Range("I66").Value = 2829.1667
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 14145.8335
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = -10713.8335
$ws.Range("N66").Value = -41864

# Row 92: columns H,J,L,N
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

# Row 100: columns H,J,L,N
$ws.Range("H100").Value = 37582.332
$ws.Range("J100").Value = 37582.332
$ws.Range("L100").Value = 37582.332
$ws.Range("N100").Value = -39746.332

# Row 102: columns H,I,K,M
$ws.Range("H102").Value = 17700.715
$ws.Range("I102").Value = 2300
$ws.Range("K102").Value = 2300
$ws.Range("M102").Value = -678

# Row 136: columns H,I,J,K,L,M,N
$ws.Range("H136").Value = 1415.7567
$ws.Range("I136").Value = 1122.7587
$ws.Range("J136").Value = 2477.875
$ws.Range("K136").Value = 3368.2761
$ws.Range("L136").Value = 7433.625
$ws.Range("M136").Value = -818.2761
$ws.Range("N136").Value = -12533.625

$ws = $wb.Worksheets.Item("BSM")
# Row 35: columns H,J,L,N
$ws.Range("H35").Value = 37756.668
$ws.Range("J35").Value = 37756.668
$ws.Range("L35").Value = 37756.668
$ws.Range("N35").Value = -38376.668

# Row 82: columns H,I,K,M
$ws.Range("H82").Value = 7590.25
$ws.Range("I82").Value = 3453.6667
$ws.Range("K82").Value = 3453.6667
$ws.Range("M82").Value = -3070.6667

# Row 85: columns H,I,K,M
$ws.Range("H85").Value = 7590.25
$ws.Range("I85").Value = 3453.6667
$ws.Range("K85").Value = 3453.6667
$ws.Range("M85").Value = -2127.6667

# Row 97: columns H,I,K,M
$ws.Range("H97").Value = 13922.8
$ws.Range("I97").Value = 4807
$ws.Range("K97").Value = 4807
$ws.Range("M97").Value = -3816

$ws = $wb.Worksheets.Item("CRP")
# Row 28: columns H,J,L,N
$ws.Range("H28").Value = 38069.5
$ws.Range("J28").Value = 38069.5
$ws.Range("L28").Value = 38069.5
$ws.Range("N28").Value = -38559.5

# Row 70: columns H,J,L,N
$ws.Range("H70").Value = 38326.6
$ws.Range("J70").Value = 38326.6
$ws.Range("L70").Value = 38326.6
$ws.Range("N70").Value = -38956.6

# Row 73: columns H,J,L,N
$ws.Range("H73").Value = 38326.6
$ws.Range("J73").Value = 38326.6
$ws.Range("L73").Value = 38326.6
$ws.Range("N73").Value = -40510.6

# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 40012.242
$ws.Range("I132").Value = 1364.5927
$ws.Range("J132").Value = 144360.9
$ws.Range("K132").Value = 4093.7781
$ws.Range("L132").Value = 433082.7
$ws.Range("M132").Value = -1563.7781
$ws.Range("N132").Value = -438142.7

# Row 134: columns H,I,J,K,L,M,N
$ws.Range("H134").Value = 702265.25
$ws.Range("I134").Value = 1823.75
$ws.Range("J134").Value = 1169226.2
$ws.Range("K134").Value = 5471.25
$ws.Range("L134").Value = 3507678.6
$ws.Range("M134").Value = -2936.25
$ws.Range("N134").Value = -3512748.6

$ws = $wb.Worksheets.Item("CUL")
# Row 15: columns H,I,J,K,L,M,N
$ws.Range("H15").Value = 130
$ws.Range("I15").Value = 120
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 360
$ws.Range("L15").Value = 420
$ws.Range("M15").Value = -220
$ws.Range("N15").Value = -700

# Row 55: columns H,I,J,K,L,M,N
$ws.Range("H55").Value = 1400
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1400
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4200
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -4554

# Row 80: columns H,I,J,K,L,M,N
$ws.Range("H80").Value = 334000320
$ws.Range("I80").Value = 2000000
$ws.Range("J80").Value = 500000500
$ws.Range("K80").Value = 6000000
$ws.Range("L80").Value = 1500001500
$ws.Range("M80").Value = -5999064
$ws.Range("N80").Value = -1500003372

# Row 83: columns H,I,J,K,L,M,N
$ws.Range("H83").Value = 334000320
$ws.Range("I83").Value = 2000000
$ws.Range("J83").Value = 500000500
$ws.Range("K83").Value = 18000000
$ws.Range("L83").Value = 4500004500
$ws.Range("M83").Value = -17995320
$ws.Range("N83").Value = -4500013860

# Row 131: columns H,I,J,K,L,M,N
$ws.Range("H131").Value = 2892.4482
$ws.Range("I131").Value = 10508.7
$ws.Range("J131").Value = 1305.7291
$ws.Range("K131").Value = 31526.1
$ws.Range("L131").Value = 3917.1873
$ws.Range("M131").Value = -26486.1
$ws.Range("N131").Value = -13997.1873

$ws = $wb.Worksheets.Item("GSM")
# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 2502.8462
$ws.Range("I132").Value = 1669.5416
$ws.Range("J132").Value = 3836.1333
$ws.Range("K132").Value = 5008.6248
$ws.Range("L132").Value = 11508.3999
$ws.Range("M132").Value = -2478.6248
$ws.Range("N132").Value = -16568.3999

$ws = $wb.Worksheets.Item("LTW")
# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 2808.4119
$ws.Range("I132").Value = 2055.8975
$ws.Range("J132").Value = 3820.4138
$ws.Range("K132").Value = 6167.6925
$ws.Range("L132").Value = 11461.2414
$ws.Range("M132").Value = -3637.6925
$ws.Range("N132").Value = -16521.2414

$ws = $wb.Worksheets.Item("WVR")
# Row 33: columns H,I,K,M
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null

# Row 36: columns H,I,K,M
$ws.Range("H36").Value = 25000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null

# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 1546.48
$ws.Range("I132").Value = 1224.5
$ws.Range("J132").Value = 3236.875
$ws.Range("K132").Value = 3673.5
$ws.Range("L132").Value = 9710.625
$ws.Range("M132").Value = -1143.5
$ws.Range("N132").Value = -14770.625

